$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "26.969.12"
$ws.Range("E2").Formula = "  +0.30%  "
$ws.Range("D3").Formula = "1.556.81"
$ws.Range("E3").Formula = "  +0.75%  "
$ws.Range("E4").Formula = "  -0.17%  "
$ws.Range("D5").Formula = "'207.26"
$ws.Range("E5").Formula = "  +0.71%  "
$ws.Range("E6").Formula = "  +0.43%  "
$ws.Range("E7").Formula = "  -0.19%  "
$ws.Range("D8").Formula = "'22.16"
$ws.Range("E8").Formula = "  +4.21%  "
$ws.Range("E9").Formula = "  +0.17%  "
$ws.Range("D10").Formula = "'0.0589"
$ws.Range("E10").Formula = "  +1.13%  "
$ws.Range("D11").Formula = "'0.0857"
$ws.Range("E11").Formula = "  -0.04%  "
$ws.Range("D12").Formula = "1.778.73"
$ws.Range("E12").Formula = "  +0.81%  "
$ws.Range("D13").Formula = "1.556.73"
$ws.Range("E13").Formula = "  +0.83%  "
$ws.Range("E14").Formula = "  +1.49%  "
$ws.Range("E15").Formula = "  +2.03%  "
$ws.Range("D16").Formula = "26.970.65"
$ws.Range("E16").Formula = "  +0.41%  "
$ws.Range("D17").Formula = "'61.81"
$ws.Range("E17").Formula = "  +0.69%  "
$ws.Range("D18").Formula = "'218.52"
$ws.Range("E18").Formula = "  +2.35%  "
$ws.Range("E19").Formula = "  +2.40%  "
$ws.Range("E20").Formula = "  +2.18%  "
$ws.Range("E21").Formula = "  -0.17%  "
$ws.Range("D22").Formula = "'4.06"
$ws.Range("E22").Formula = "  +1.26%  "
$ws.Range("E23").Formula = "  +0.84%  "
$ws.Range("E24").Formula = "  +0.53%  "
$ws.Range("D25").Formula = "'154.23"
$ws.Range("E25").Formula = "  +1.09%  "
$ws.Range("E26").Formula = "  +0.95%  "
$ws.Range("D27").Formula = "'14.99"
$ws.Range("E27").Formula = "  +1.29%  "
$ws.Range("E28").Formula = "  +1.32%  "
$ws.Range("E29").Formula = "  -0.13%  "
$ws.Range("E30").Formula = "  +2.67%  "
$ws.Range("E31").Formula = "  -0.14%  "
$ws.Range("E32").Formula = "  +0.95%  "
$ws.Range("D33").Formula = "1.426.07"
$ws.Range("E33").Formula = "  +5.40%  "
$ws.Range("D34").Formula = "'3.09"
$ws.Range("E34").Formula = "  +5.29%  "
$ws.Range("D35").Formula = "'1.59"
$ws.Range("E35").Formula = "  +4.10%  "
$ws.Range("E36").Formula = "  +1.82%  "
$ws.Range("E37").Formula = "  +0.27%  "
$ws.Range("E38").Formula = "  +1.00%  "
$ws.Range("D39").Formula = "'0.523"
$ws.Range("E39").Formula = "  +0.50%  "
$ws.Range("E40").Formula = "  +1.13%  "
$ws.Range("D41").Formula = "'5.76"
$ws.Range("E41").Formula = "  +3.43%  "
$ws.Range("E42").Formula = "  -0.15%  "
$ws.Range("E43").Formula = "  +4.86%  "
$ws.Range("E44").Formula = "  -0.20%  "
$ws.Range("D45").Formula = "'64.62"
$ws.Range("E45").Formula = "  +2.06%  "
$ws.Range("D46").Formula = "'1.77"
$ws.Range("E46").Formula = "  +2.80%  "
$ws.Range("D47").Formula = "1.692.29"
$ws.Range("E47").Formula = "  +0.81%  "
$ws.Range("D48").Formula = "'88.04"
$ws.Range("E48").Formula = "  +2.62%  "
$ws.Range("D49").Formula = "'0.0521"
$ws.Range("E49").Formula = "  +2.09%  "
$ws.Range("D50").Formula = "0.0₇0996"
$ws.Range("E50").Formula = "  +2.31%  "
$ws.Range("D51").Formula = "'0.0957"
$ws.Range("E51").Formula = "  +1.06%  "
